$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 9454
$ws.Range("F5").Value = 733
$ws.Range("F6").Value = 617
$ws.Range("F7").Value = 208
$ws.Range("F8").Value = 322
$ws.Range("F11").Value = 1621
$ws.Range("F12").Value = 1402
$ws.Range("F15").Value = 1465
$ws.Range("F16").Value = 122
$ws.Range("F17").Value = 299
$ws.Range("F19").Value = 135
$ws.Range("F21").Value = 375
$ws.Range("F22").Value = 1107
$ws.Range("F25").Value = 47
$ws.Range("F26").Value = 275
$ws.Range("F28").Value = 253
$ws.Range("F31").Value = 631
$ws.Range("F32").Value = 3
$ws.Range("F35").Value = 177
$ws.Range("F36").Value = 315
$ws.Range("F38").Value = 228
$ws.Range("F39").Value = 603
$ws.Range("F40").Value = 511
$ws.Range("F42").Value = 733

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 36
$ws.Range("F6").Value = 66
$ws.Range("F10").Value = 686
$ws.Range("F12").Value = 46
$ws.Range("F19").Value = 30
$ws.Range("F20").Value = 1059
$ws.Range("F21").Value = 279
$ws.Range("F23").Value = 8
$ws.Range("F25").Value = 304
$ws.Range("F29").Value = 19
$ws.Range("F37").Value = 15
$ws.Range("F38").Value = 23

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 769
$ws.Range("F5").Value = 361
$ws.Range("F6").Value = 157
$ws.Range("F7").Value = 2353
$ws.Range("F8").Value = 3564
$ws.Range("F9").Value = 5
$ws.Range("F11").Value = 60
$ws.Range("F12").Value = 98

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 9454
$ws.Range("F4").Value = 361
$ws.Range("F5").Value = 157
$ws.Range("F6").Value = 3564
$ws.Range("F7").Value = 733
$ws.Range("F8").Value = 60
$ws.Range("F9").Value = 60
$ws.Range("F10").Value = 617
$ws.Range("F11").Value = 208
$ws.Range("F12").Value = 322
$ws.Range("F13").Value = 686
$ws.Range("F14").Value = 1402
$ws.Range("F16").Value = 98
$ws.Range("F17").Value = 98
$ws.Range("F18").Value = 1465
$ws.Range("F19").Value = 299
$ws.Range("F21").Value = 135
$ws.Range("F22").Value = 1107
$ws.Range("F24").Value = 46
$ws.Range("F27").Value = 47
$ws.Range("F28").Value = 275
$ws.Range("F29").Value = 30
$ws.Range("F31").Value = 253
$ws.Range("F32").Value = 1059
$ws.Range("F33").Value = 279
$ws.Range("F35").Value = 631
$ws.Range("F36").Value = 8
$ws.Range("F37").Value = 3
$ws.Range("F39").Value = 304
$ws.Range("F40").Value = 304
$ws.Range("F41").Value = 315
$ws.Range("F43").Value = 228
$ws.Range("F45").Value = 603
$ws.Range("F46").Value = 511
$ws.Range("F47").Value = 733
$ws.Range("F52").Value = 15
